$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value  = 5.344
$ws.Range("A9").Value  = -21.723
$ws.Range("B12").Value = 5.57
$ws.Range("A18").Value = -22.156
$ws.Range("A20").Value = -20.391
$ws.Range("B26").Value = 5.761
$ws.Range("A27").Value = -21.188
$ws.Range("B27").Value = 5.695000000000001
$ws.Range("B29").Value = 5.688
$ws.Range("B37").Value = 8.376000000000001
$ws.Range("B38").Value = 4.947
$ws.Range("B51").Value = 5.82
$ws.Range("B55").Value = 5.867
$ws.Range("A69").Value = -21.61
$ws.Range("B69").Value = 5.688000000000001
$ws.Range("B70").Value = 5.233
$ws.Range("A76").Value = -20.396
$ws.Range("A82").Value = -22.152
$ws.Range("B83").Value = 5.695
$ws.Range("B102").Value = 7.231
